$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell "time_taken" in column F, styled like the existing
# header row (bold + border), by copying the format from E1.
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Timestamp values for each data row (2-25)
$timestamps = @(
    "2021-10-05 13:38:47.063245",
    "2021-10-05 13:38:47.063253",
    "2021-10-05 13:38:47.063256",
    "2021-10-05 13:38:47.063258",
    "2021-10-05 13:38:47.063260",
    "2021-10-05 13:38:47.063262",
    "2021-10-05 13:38:47.063264",
    "2021-10-05 13:38:47.063266",
    "2021-10-05 13:38:47.063268",
    "2021-10-05 13:38:47.063270",
    "2021-10-05 13:38:47.063272",
    "2021-10-05 13:38:47.063274",
    "2021-10-05 13:38:47.063276",
    "2021-10-05 13:38:47.063278",
    "2021-10-05 13:38:47.063280",
    "2021-10-05 13:38:47.063282",
    "2021-10-05 13:38:47.063284",
    "2021-10-05 13:38:47.063286",
    "2021-10-05 13:38:47.063288",
    "2021-10-05 13:38:47.063290",
    "2021-10-05 13:38:47.063292",
    "2021-10-05 13:38:47.063293",
    "2021-10-05 13:38:47.063295",
    "2021-10-05 13:38:47.063297"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
